$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - new purchase entry: 3D printing & Laser Cutting
$ws.Range("B10").Value = 43537
$ws.Range("C10").Value = "3D printing & Laser Cutting"
$ws.Range("E10").Value = "F.Holmes "
$ws.Range("D10").Value = "Internal"
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 1
$ws.Range("H10").Formula = "=F10"
$ws.Range("I10").Value = "Received"
$ws.Range("J10").Formula = "=IF(G10=0,`"`",IF(I10=`"Received`",J9+H10,J9))"
$ws.Range("L10").Value = "Tri Track and Sorting Rig"

# Row 11 - new purchase entry: Kinect
$ws.Range("B11").Value = 43538
$ws.Range("C11").Value = "Kinect"
$ws.Range("E11").Value = "F.Holmes "
$ws.Range("D11").Value = "Internal"
$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 1
$ws.Range("H11").Formula = "=F11"
$ws.Range("I11").Value = "Received"
$ws.Range("J11").Formula = "=IF(G11=0,`"`",IF(I11=`"Received`",J10+H11,J10))"
$ws.Range("L11").Value = "Tri Track"

# Leave the cursor where the author last left it
$null = $ws.Range("L19").Select()

